$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# WA (workaround) for "convert arguments from code to visual mode" bug report:
# three new test rows appended after the existing table (rows 44-46),
# plus the remaining numbering rows (47-56) that were already blank below.

$code = "def func0(arg0):`n  if compare(2, `"!=`", arg0):`n    if checkFilled(1, `">`", 2):`n      sub(1, 2, arg0)`n    else:`n      sub(2, 1, arg0)`ndef sub(arg0, arg1, arg2):`n  pour(arg0, arg1)`n  if compare(arg1, `"!=`", arg2):`n    pourOut(arg1)`nfunc0(5)"

# Row 44
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "применять изменения имен функций и аргументов на Enter"
$ws.Range("B44").WrapText = $true
$ws.Range("I44").Value = "NG"

# Row 45
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "Некорректный аргумент`n" + $code
$ws.Range("B45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 16.5
$ws.Range("I45").Value = "OK"

# Row 46
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "не конвертируется в команды`n" + $code + "`n"
$ws.Range("B46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 15
$ws.Range("I46").Value = "NG"

# Rows 47-56: just the running counter in column A, no other content.
For ($i = 47; $i -le 56; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

$ws.Range("I46").Select()
